$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 14108.469
$ws.Range("I15").Value = 14108.469
$ws.Range("K15").Value = 42325.407
$ws.Range("M15").Value = -42156.407
$ws.Range("H18").Value = 576.6667
$ws.Range("I18").Value = 576.6667
$ws.Range("K18").Value = 576.6667
$ws.Range("M18").Value = -292.6667
$ws.Range("H19").Value = 273.45456
$ws.Range("I19").Value = 237.81818
$ws.Range("K19").Value = 237.81818
$ws.Range("M19").Value = -62.81818000000001
$ws.Range("H29").Value = 16707.143
$ws.Range("I29").Value = 2000
$ws.Range("J29").Value = 22590
$ws.Range("K29").Value = 6000
$ws.Range("L29").Value = 67770
$ws.Range("M29").Value = -5719
$ws.Range("N29").Value = -68332
$ws.Range("H55").Value = 350.0909
$ws.Range("I55").Value = 516.8333
$ws.Range("J55").Value = 150
$ws.Range("K55").Value = 516.8333
$ws.Range("L55").Value = 150
$ws.Range("M55").Value = -302.8333
$ws.Range("N55").Value = -578
$ws.Range("H125").Value = 4810.125
$ws.Range("I125").Value = 1612.8
$ws.Range("J125").Value = 6263.4546
$ws.Range("K125").Value = 14515.2
$ws.Range("L125").Value = 56371.0914
$ws.Range("M125").Value = -12055.2
$ws.Range("N125").Value = -61291.0914
$ws.Range("H132").Value = 1287.25
$ws.Range("I132").Value = 723.2381
$ws.Range("K132").Value = 2169.7143
$ws.Range("M132").Value = 360.2856999999999
$ws.Range("H133").Value = 29000
$ws.Range("J133").Value = 29000
$ws.Range("L133").Value = 29000
$ws.Range("N133").Value = -39120
$ws.Range("H134").Value = 34111.11
$ws.Range("J134").Value = 34111.11
$ws.Range("L134").Value = 34111.11
$ws.Range("N134").Value = -44251.11
$ws.Range("H135").Value = 2580.7693
$ws.Range("I135").Value = 2783.0908
$ws.Range("J135").Value = 1468
$ws.Range("K135").Value = 25047.8172
$ws.Range("L135").Value = 13212
$ws.Range("M135").Value = -22512.8172
$ws.Range("N135").Value = -18282
$ws.Range("H136").Value = 50000
$ws.Range("J136").Value = 50000
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200
$ws.Range("H137").Value = 1789.5264
$ws.Range("I137").Value = 1721.5
$ws.Range("J137").Value = 1980
$ws.Range("K137").Value = 5164.5
$ws.Range("L137").Value = 5940
$ws.Range("M137").Value = -2614.5
$ws.Range("N137").Value = -11040
$ws.Range("H138").Value = 3259.6667
$ws.Range("I138").Value = 1163.4667
$ws.Range("J138").Value = 5006.5
$ws.Range("K138").Value = 3490.4001
$ws.Range("L138").Value = 15019.5
$ws.Range("M138").Value = 1649.5999
$ws.Range("N138").Value = -25299.5
$ws.Range("H139").Value = 26640
$ws.Range("J139").Value = 26640
$ws.Range("L139").Value = 26640
$ws.Range("N139").Value = -36920
$ws.Range("H140").Value = 48750
$ws.Range("J140").Value = 48750
$ws.Range("L140").Value = 48750
$ws.Range("N140").Value = -59110
$ws.Range("H141").Value = 2536.4285
$ws.Range("I141").Value = 2181.389
$ws.Range("K141").Value = 6544.167
$ws.Range("M141").Value = -1364.167

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 5004
$ws.Range("I19").Value = 8
$ws.Range("K19").Value = 8
$ws.Range("M19").Value = 221
$ws.Range("H33").Value = 4513
$ws.Range("I33").Value = 26
$ws.Range("K33").Value = 26
$ws.Range("M33").Value = 303
$ws.Range("H36").Value = 1666.6666
$ws.Range("I36").Value = 1666.6666
$ws.Range("K36").Value = 1666.6666
$ws.Range("M36").Value = -1320.6666
$ws.Range("H44").Value = 18100
$ws.Range("J44").Value = 18100
$ws.Range("L44").Value = 18100
$ws.Range("N44").Value = -19076
$ws.Range("H55").Value = 11718.2
$ws.Range("J55").Value = 21795.5
$ws.Range("L55").Value = 21795.5
$ws.Range("N55").Value = -22425.5
$ws.Range("H97").Value = 1311.5
$ws.Range("I97").Value = 967.5
$ws.Range("K97").Value = 967.5
$ws.Range("M97").Value = -471.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 2237.8
$ws.Range("I25").Value = 2237.8
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 2237.8
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -2002.8
$ws.Range("N25").ClearContents()
$ws.Range("H55").Value = 30779
$ws.Range("J55").Value = 30779
$ws.Range("L55").Value = 30779
$ws.Range("N55").Value = -31325
$ws.Range("H94").Value = 1880.75
$ws.Range("I94").Value = 1562.8334
$ws.Range("J94").Value = 2198.6667
$ws.Range("K94").Value = 1562.8334
$ws.Range("L94").Value = 2198.6667
$ws.Range("M94").Value = -1111.8334
$ws.Range("N94").Value = -3100.6667
$ws.Range("H134").Value = 4859.75
$ws.Range("I134").Value = 5906.9614
$ws.Range("J134").Value = 2914.9285
$ws.Range("K134").Value = 17720.8842
$ws.Range("L134").Value = 8744.7855
$ws.Range("M134").Value = -15185.8842
$ws.Range("N134").Value = -13814.7855

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 50668.668
$ws.Range("J3").Value = 50668.668
$ws.Range("L3").Value = 50668.668
$ws.Range("N3").Value = -50894.668
$ws.Range("H25").Value = 1000000000
$ws.Range("J25").Value = 1000000000
$ws.Range("L25").Value = 1000000000
$ws.Range("N25").Value = -1000000348

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 97.42856999999999
$ws.Range("I6").Value = 97.42856999999999
$ws.Range("K6").Value = 292.28571
$ws.Range("M6").Value = -179.28571
$ws.Range("H107").Value = 392.7143
$ws.Range("J107").Value = 510.8889
$ws.Range("L107").Value = 1532.6667
$ws.Range("N107").Value = -5372.6667
$ws.Range("H122").Value = 931.3333
$ws.Range("I122").Value = 747.1667
$ws.Range("J122").Value = 1299.6666
$ws.Range("K122").Value = 6724.5003
$ws.Range("L122").Value = 11696.9994
$ws.Range("M122").Value = -4274.5003
$ws.Range("N122").Value = -16596.9994
$ws.Range("H124").Value = 6899.778
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 6899.778
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 20699.334
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -30519.334
$ws.Range("H132").Value = 1739143.2
$ws.Range("I132").Value = 2500
$ws.Range("J132").Value = 1795164.1
$ws.Range("K132").Value = 22500
$ws.Range("L132").Value = 16156476.9
$ws.Range("M132").Value = -19970
$ws.Range("N132").Value = -16161536.9
$ws.Range("H136").Value = 5052.108
$ws.Range("I136").Value = 50515
$ws.Range("K136").Value = 151545
$ws.Range("M136").Value = -146445
$ws.Range("H137").Value = 13267.733
$ws.Range("I137").Value = 27272.25
$ws.Range("J137").Value = 8175.1816
$ws.Range("K137").Value = 81816.75
$ws.Range("L137").Value = 24525.5448
$ws.Range("M137").Value = -76716.75
$ws.Range("N137").Value = -34725.5448
$ws.Range("H138").Value = 18356.428
$ws.Range("I138").Value = 21077
$ws.Range("J138").Value = 2033
$ws.Range("K138").Value = 63231
$ws.Range("L138").Value = 6099
$ws.Range("M138").Value = -58091
$ws.Range("N138").Value = -16379
$ws.Range("H139").Value = 4687.1396
$ws.Range("I139").Value = 8539.071
$ws.Range("J139").Value = 2827.5862
$ws.Range("K139").Value = 25617.213
$ws.Range("L139").Value = 8482.758600000001
$ws.Range("M139").Value = -20477.213
$ws.Range("N139").Value = -18762.7586
$ws.Range("H140").Value = 2467.6
$ws.Range("I140").Value = 2505.7856
$ws.Range("J140").Value = 1933
$ws.Range("K140").Value = 7517.3568
$ws.Range("L140").Value = 5799
$ws.Range("M140").Value = -2337.3568
$ws.Range("N140").Value = -16159
$ws.Range("H141").Value = 18460.1
$ws.Range("I141").Value = 31700.25
$ws.Range("K141").Value = 95100.75
$ws.Range("M141").Value = -89920.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 500030000
$ws.Range("I29").Value = 500030000
$ws.Range("K29").Value = 500030000
$ws.Range("M29").Value = -500029710

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2360
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 2825
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 2825
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -3415
$ws.Range("H27").Value = 2360
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 2825
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 2825
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -3039
$ws.Range("H46").Value = 1055.5555
$ws.Range("I46").Value = 1037.5
$ws.Range("J46").Value = 1070
$ws.Range("K46").Value = 1037.5
$ws.Range("L46").Value = 1070
$ws.Range("M46").Value = -849.5
$ws.Range("N46").Value = -1446
$ws.Range("H55").Value = 467.33334
$ws.Range("I55").Value = 149.75
$ws.Range("J55").Value = 626.125
$ws.Range("K55").Value = 149.75
$ws.Range("L55").Value = 626.125
$ws.Range("M55").Value = 23.25
$ws.Range("N55").Value = -972.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 4500
$ws.Range("I32").Value = 3500
$ws.Range("J32").Value = 6500
$ws.Range("K32").Value = 3500
$ws.Range("L32").Value = 6500
$ws.Range("M32").Value = -3183
$ws.Range("N32").Value = -7134
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H123").Value = 36500
$ws.Range("J123").Value = 43000
$ws.Range("L123").Value = 43000
$ws.Range("N123").Value = -52800
$ws.Range("H136").Value = 1102.5862
$ws.Range("I136").Value = 698.63635
$ws.Range("J136").Value = 2372.1428
$ws.Range("K136").Value = 2095.90905
$ws.Range("L136").Value = 7116.428400000001
$ws.Range("M136").Value = 454.0909499999998
$ws.Range("N136").Value = -12216.4284
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 58857.25
$ws.Range("J140").Value = 58857.25
$ws.Range("L140").Value = 58857.25
$ws.Range("N140").Value = -69217.25
